$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the two new Cypher query strings (multi-line) ---
$casesQuery = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.race = `"BLACK_OR_AFRICAN_AMERICAN`"`nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

$statQuery = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.race = `"BLACK_OR_AFRICAN_AMERICAN`"`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

# --- Insert a new column A in front of everything, shifting the old A:D to B:E ---
$ws.Columns.Item(1).Insert()

# --- New column A header/value ---
$ws.Cells.Item(1, 1).Value = "TabName"
$ws.Cells.Item(2, 1).Value = "CasesTab"

# --- Replace the two query cells in place with the new Cypher text ---
$ws.Cells.Item(2, 2).Value = $casesQuery
$ws.Cells.Item(2, 3).Value = $statQuery

# --- Column A width (narrow, best-fit to "TabName"/"CasesTab") ---
$ws.Columns.Item(1).ColumnWidth = 8.81640625

# --- Row height: re-autofit row 2 to account for the taller wrapped text ---
$ws.Rows.Item(2).AutoFit()

# --- Selection, matching the saved cursor position in the edited file ---
$ws.Range("B5").Select()
